$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Row 16 ("GLAY") is written FIRST, before C2/C3 are touched, so
#    the shared-string table reuses the slot "SEVENTEEN LIVE"
#    currently occupies once that text is freed below.
# ------------------------------------------------------------------

# Reusable "icon" source cells (already-typed emoji, avoids re-typing
# the characters and keeps the shared-string identity intact).
$redIcon = $ws.Range("B2")
$blueIcon = $ws.Range("B6")

function Add-EventRow($Row, $ExcelDate, $IconSource, $Name, $Wrap) {
    $dateCell = $ws.Cells.Item($Row, 1)
    $ws.Range("A2").Copy()
    $dateCell.PasteSpecial(-4122)
    $dateCell.Value = $ExcelDate

    $iconCell = $ws.Cells.Item($Row, 2)
    $IconSource.Copy()
    $iconCell.PasteSpecial(-4122)
    $iconCell.Value = $IconSource.Value2

    $nameCell = $ws.Cells.Item($Row, 3)
    $nameCell.Value = $Name
    if ($Wrap) {
        $nameCell.WrapText = $true
    }
}

Add-EventRow 16 45816 $redIcon  "GLAY"         $false

# ------------------------------------------------------------------
# 2) Repoint the two existing "SEVENTEEN LIVE" cells at the new
#    "SEVENTEEN" text right away (immediately after "GLAY" is minted
#    above, before any other new unique string is introduced) - this
#    frees up the old shared-string slot for "GLAY" and creates the
#    fresh "SEVENTEEN" slot directly after it, matching the target
#    shared-string ordering.
# ------------------------------------------------------------------
$ws.Range("C2").Value = "SEVENTEEN"
$ws.Range("C3").Value = "SEVENTEEN"

Add-EventRow 17 45840 $redIcon  "三代目JSOUL"   $false
Add-EventRow 18 45841 $redIcon  "三代目JSOUL"   $false
Add-EventRow 19 45850 $redIcon  "JIN（BTS)"     $false
Add-EventRow 20 45851 $redIcon  "JIN（BTS)"     $false
Add-EventRow 21 45857 $redIcon  "king&prince"  $false
Add-EventRow 22 45858 $redIcon  "king&prince"  $false
Add-EventRow 23 45859 $redIcon  "king&prince"  $false
Add-EventRow 24 45871 $blueIcon "ENHYPEN"      $true
Add-EventRow 25 45872 $blueIcon "ENHYPEN"      $true
Add-EventRow 26 45892 $redIcon  "櫻坂46"        $false
Add-EventRow 27 45893 $redIcon  "櫻坂47"        $false
Add-EventRow 28 45899 $redIcon  "福山雅治"      $false
Add-EventRow 29 45900 $redIcon  "福山雅治"      $false
Add-EventRow 30 45913 $blueIcon "ワンオクロック" $false
Add-EventRow 31 45914 $blueIcon "ワンオクロック" $false
Add-EventRow 32 45921 $redIcon  "星野源"        $false
Add-EventRow 33 46026 $redIcon  "Saucy Dog"    $false
Add-EventRow 34 46095 $redIcon  "Ｖａｕｎｄｙ"     $false

# ------------------------------------------------------------------
# 3) Selection / view bookkeeping to match the saved workbook state.
# ------------------------------------------------------------------
$ws.Range("C7").Select()
